# Update the localization-status report: the files
#   0c92f6c4-7ab7-4dc9-aa23-664de468b52f.md   (row 2)
#   1aa451c0-4d8f-415b-9f13-736fa576ed5a.md   (row 3)
# moved from "Ready for handoff" to "In Translation".
# Reflect the new status on the per-locale sheets (zh-cn, de-de) and on
# the rolled-up Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "In Translation"
$overview.Range("C2").Value = "In Translation"
$overview.Range("B3").Value = "In Translation"
$overview.Range("C3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "In Translation"
$zhcn.Range("B3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "In Translation"
$dede.Range("B3").Value = "In Translation"
